# Update market/profit figures on the Yojimbo profits workbook.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 26: Night Squawker / Iron Lantern Shield
$ws.Range("H26").Value = 659.8
$ws.Range("I26").Value = 659.8
$ws.Range("K26").Value = 659.8
$ws.Range("M26").Value = -329.8

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 3138.1
$ws.Range("I74").Value = 2929.6667
$ws.Range("J74").Value = 5014
$ws.Range("K74").Value = 2929.6667
$ws.Range("L74").Value = 5014
$ws.Range("M74").Value = -2055.6667
$ws.Range("N74").Value = -6762

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 3138.1
$ws.Range("I77").Value = 2929.6667
$ws.Range("J77").Value = 5014
$ws.Range("K77").Value = 14648.3335
$ws.Range("L77").Value = 25070
$ws.Range("M77").Value = -10280.3335
$ws.Range("N77").Value = -33806

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1408.8206
$ws.Range("I132").Value = 1025.1428
$ws.Range("J132").Value = 2385.4546
$ws.Range("K132").Value = 3075.4284
$ws.Range("L132").Value = 7156.3638
$ws.Range("M132").Value = -545.4284000000002
$ws.Range("N132").Value = -12216.3638

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1104.1923
$ws.Range("I99").Value = 996.2083
$ws.Range("K99").Value = 996.2083
$ws.Range("M99").Value = 501.7917

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1180.2667
$ws.Range("I134").Value = 952.8
$ws.Range("K134").Value = 2858.4
$ws.Range("M134").Value = -323.3999999999996

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1235.7858
$ws.Range("I16").Value = 1291.1
$ws.Range("J16").Value = 1097.5
$ws.Range("K16").Value = 1291.1
$ws.Range("L16").Value = 1097.5
$ws.Range("M16").Value = -1004.1
$ws.Range("N16").Value = -1671.5

# Row 26: As the Worm Turns / Yew Radical
$ws.Range("H26").Value = 22500
$ws.Range("J26").Value = 22500
$ws.Range("L26").Value = 22500
$ws.Range("N26").Value = -23074

# Row 29: Grinding It Out / Mudstone Grinding Wheel
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1666
$ws.Range("I107").Value = 1785.5294
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 1785.5294
$ws.Range("L107").Value = 650
$ws.Range("M107").Value = 134.4706000000001
$ws.Range("N107").Value = -4490

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1235.7858
$ws.Range("I113").Value = 1291.1
$ws.Range("J113").Value = 1097.5
$ws.Range("K113").Value = 1291.1
$ws.Range("L113").Value = 1097.5
$ws.Range("M113").Value = 878.9000000000001
$ws.Range("N113").Value = -5437.5

# Row 114: Ground to a Halt / White Ash Grinding Wheel
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 988.8333
$ws.Range("I132").Value = 652.43243
$ws.Range("K132").Value = 1957.29729
$ws.Range("M132").Value = 572.70271

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1736.25
$ws.Range("I134").Value = 1301.3158
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 3903.9474
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -1368.9474
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CUL")
# Row 36: Love's Crumpets Lost / Crumpet
$ws.Range("H36").Value = 4920.5713
$ws.Range("I36").Value = 740.6667
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 2222.0001
$ws.Range("L36").Value = 90000
$ws.Range("M36").Value = -2053.0001
$ws.Range("N36").Value = -90338

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1221795.1
$ws.Range("I131").Value = 379
$ws.Range("J131").Value = 1452251
$ws.Range("K131").Value = 1137
$ws.Range("L131").Value = 4356753
$ws.Range("M131").Value = 3903
$ws.Range("N131").Value = -4366833

$ws = $wb.Worksheets.Item("GSM")
# Row 24: Bad Guys Eat Brass / Brass Ring of Crafting
$ws.Range("H24").Value = 60007
$ws.Range("J24").Value = 60007
$ws.Range("L24").Value = 60007
$ws.Range("N24").Value = -60353

# Row 29: Music to Their Ears / Brass Ear Cuffs
$ws.Range("H29").Value = 27499.5
$ws.Range("J29").Value = 27499.5
$ws.Range("L29").Value = 27499.5
$ws.Range("N29").Value = -28079.5

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("I113").Value = 793.8182
$ws.Range("J113").Value = 775
$ws.Range("K113").Value = 793.8182
$ws.Range("L113").Value = 775
$ws.Range("M113").Value = 1376.1818
$ws.Range("N113").Value = -5115

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2751.75
$ws.Range("I122").Value = 2503.5
$ws.Range("K122").Value = 7510.5
$ws.Range("M122").Value = -5060.5

# Row 134: Guaranteed Gem / Ihuykanite
$ws.Range("H134").Value = 5945.6
$ws.Range("J134").Value = 5945.6
$ws.Range("L134").Value = 17836.8
$ws.Range("N134").Value = -22906.8

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 846866.4
$ws.Range("I22").Value = 143608.72
$ws.Range("J22").Value = 1667333.6
$ws.Range("K22").Value = 143608.72
$ws.Range("L22").Value = 1667333.6
$ws.Range("M22").Value = -143313.72
$ws.Range("N22").Value = -1667923.6

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 846866.4
$ws.Range("I27").Value = 143608.72
$ws.Range("J27").Value = 1667333.6
$ws.Range("K27").Value = 143608.72
$ws.Range("L27").Value = 1667333.6
$ws.Range("M27").Value = -143501.72
$ws.Range("N27").Value = -1667547.6

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 345.95
$ws.Range("I55").Value = 258.64285
$ws.Range("J55").Value = 549.6667
$ws.Range("K55").Value = 258.64285
$ws.Range("L55").Value = 549.6667
$ws.Range("M55").Value = -85.64285000000001
$ws.Range("N55").Value = -895.6667

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1814.2858
$ws.Range("I61").Value = 1560
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 1560
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -1358
$ws.Range("N61").Value = -2854

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1814.2858
$ws.Range("I113").Value = 1560
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 1560
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = 610
$ws.Range("N113").Value = -6790

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6145.7144
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 3836.6667
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 11510.0001
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -16410.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1410.0845
$ws.Range("I132").Value = 1163.8518
$ws.Range("J132").Value = 2192.2354
$ws.Range("K132").Value = 3491.5554
$ws.Range("L132").Value = 6576.706200000001
$ws.Range("M132").Value = -961.5553999999997
$ws.Range("N132").Value = -11636.7062

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1013.1923
$ws.Range("I136").Value = 1054.3334
$ws.Range("J136").Value = 920.625
$ws.Range("K136").Value = 3163.0002
$ws.Range("L136").Value = 2761.875
$ws.Range("M136").Value = -613.0001999999999
$ws.Range("N136").Value = -7861.875
